$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the values currently in column A (A1:A7) into column C (C1:C7).
for ($r = 1; $r -le 7; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 3).Value2 = $val
}

# Row 7 also gets duplicated into column B.
$ws.Cells.Item(7, 2).Value2 = $ws.Cells.Item(7, 1).Value2

# Clear the old header cell in column A (it now only lives in C1).
$ws.Cells.Item(1, 1).ClearContents()

# Update the active selection to match the final state (C7).
$ws.Range("C7").Select()
